$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bug fix: the "data_saved_notify" mail body (row 8, column F) referenced a
# non-existent placeholder ${updated_user} where it should say ${target_user}
# (matching the other ${target_user} placeholders used later in the same text).
$ws.Range("F8").Value = 'Data for ${target_table} has been ${create_or_update} by user ${target_user}.\nPlease check the following.\n\n${create_or_update} User: ${target_user}\n${create_or_update} Date: ${target_datetime}\n${create_or_update} Data: ${value_url/link=true}\n${free_space}'

# The author's cursor ended up on F9 when the file was saved.
$ws.Range("F9").Select()
